$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.437828183174133
$ws.Range("B1").Value = 2.002362728118896
$ws.Range("C1").Value = 3.559284925460815
$ws.Range("D1").Value = 3.522897005081177
$ws.Range("E1").Value = 0.8067057132720947
